$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the active cell selection on the sheet view
$ws.Range("F3").Select()

# Clear computed values in column G for rows 9-12 and 15-20 (keep styles)
$ws.Range("G9").ClearContents()
$ws.Range("G10").ClearContents()
$ws.Range("G11").ClearContents()
$ws.Range("G12").ClearContents()
$ws.Range("G15").ClearContents()
$ws.Range("G16").ClearContents()
$ws.Range("G17").ClearContents()
$ws.Range("G18").ClearContents()
$ws.Range("G19").ClearContents()
$ws.Range("G20").ClearContents()

# Clear H17 and H18 values
$ws.Range("H17").ClearContents()
$ws.Range("H18").ClearContents()

# Clear F22, F23, F24 values
$ws.Range("F22").ClearContents()
$ws.Range("F23").ClearContents()
$ws.Range("F24").ClearContents()
